# Updates the 江西-漫展信息 workbook:
#   - bumps a handful of "想去人数" (interest count) figures on both the
#     "展览" and "全部类型" sheets
#   - renames 诚瑞橙子运动馆 -> 诚瑞橙子体育馆 in the location text
#   - appends a new convention row (赣州·第二届异次元动漫嘉年华) at the
#     bottom of both sheets

$wb = $excel.ActiveWorkbook

function Set-WantCount {
    param($ws, [int]$row, [double]$newValue)
    $ws.Cells.Item($row, 6).Value = $newValue
}

function Set-NewRow {
    param($ws, [int]$row)

    # Copy the bold/centered/bordered style used by every "index" cell in
    # column A (style index 1 in styles.xml) from the row above instead of
    # re-describing it by hand, so no duplicate style entry is created.
    $ws.Range($ws.Cells.Item($row - 1, 1), $ws.Cells.Item($row - 1, 1)).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $row - 1

    # B/date-like text: forcing a text number-format keeps Excel from
    # auto-converting the "2024-08-08" literal into a date serial, then
    # resetting the style back to Normal keeps the cell style identical
    # to its untouched neighbours (no stray "s" attribute left behind).
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = "2024-08-08"
    $bCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "赣州·第二届异次元动漫嘉年华"
    $ws.Cells.Item($row, 4).Value = "金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆"
    $ws.Cells.Item($row, 5).Value = "2024.08.08 08:00-08.08 17:00"
    $ws.Cells.Item($row, 6).Value = 1
    $ws.Cells.Item($row, 7).Value = 45
    $ws.Cells.Item($row, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84184"
    $ws.Cells.Item($row, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/F5F9vvqX1712758945373.jpeg"
}

# ---- "展览" sheet (rows 1-33 before the edit) ----
$wsExpo = $wb.Worksheets.Item("展览")

Set-WantCount $wsExpo 3 192
Set-WantCount $wsExpo 5 5159
Set-WantCount $wsExpo 9 573
Set-WantCount $wsExpo 10 528
Set-WantCount $wsExpo 12 20
Set-WantCount $wsExpo 14 4110
Set-WantCount $wsExpo 18 91
Set-WantCount $wsExpo 19 3104
Set-WantCount $wsExpo 20 151
Set-WantCount $wsExpo 21 1049
Set-WantCount $wsExpo 25 90
Set-WantCount $wsExpo 26 23
Set-WantCount $wsExpo 28 67

$wsExpo.Cells.Item(30, 4).Value = "迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子体育馆"
Set-WantCount $wsExpo 30 15
Set-WantCount $wsExpo 32 9
Set-WantCount $wsExpo 33 9

Set-NewRow $wsExpo 34

# ---- "全部类型" sheet (rows 1-34 before the edit, offset by +1 row vs 展览) ----
$wsAll = $wb.Worksheets.Item("全部类型")

Set-WantCount $wsAll 3 192
Set-WantCount $wsAll 6 5159
Set-WantCount $wsAll 10 573
Set-WantCount $wsAll 11 528
Set-WantCount $wsAll 13 20
Set-WantCount $wsAll 15 4110
Set-WantCount $wsAll 19 91
Set-WantCount $wsAll 20 3104
Set-WantCount $wsAll 21 151
Set-WantCount $wsAll 22 1049
Set-WantCount $wsAll 26 90
Set-WantCount $wsAll 27 23
Set-WantCount $wsAll 29 67

$wsAll.Cells.Item(31, 4).Value = "迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子体育馆"
Set-WantCount $wsAll 31 15
Set-WantCount $wsAll 33 9
Set-WantCount $wsAll 34 9

Set-NewRow $wsAll 35
